$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E edits on rows 2-25 (row numbers unaffected by later row deletions) ---
$ws.Range("E2").Value = -7.2
$ws.Range("E6").ClearContents()
$ws.Range("E12").Value = -5.3
$ws.Range("E14").ClearContents()
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E23").ClearContents()
$ws.Range("E24").ClearContents()

# --- Remove row 26 ("RM 232") entirely ---
$ws.Rows.Item(26).Delete()

# --- Remove the row that is now "SC 92" (originally row 28, now row 27) ---
$ws.Rows.Item(27).Delete()

# --- After the two deletions, fix up B/E values on the shifted rows ---
# New row 26 = SC 5
$ws.Range("B26").Value = -20.2
# New row 27 = SC 101
$ws.Range("B27").ClearContents()
# New row 30 = SC 120
$ws.Range("B30").Value = -19.7
# New row 31 = SC 132
$ws.Range("E31").Value = -8.1
# New row 32 = SC 193
$ws.Range("B32").ClearContents()
# New row 33 = SC 232
$ws.Range("E33").Value = -10.7

# --- Update the sheet dimension to match the new extent ---
$ws.Range("A1:F33").Value = $ws.Range("A1:F33").Value()
